# Insert two new data rows (48 and 49) for variety "Artic Pride" above the
# existing row 48 ("Nectar Crest"), pushing the former rows 48-138 down to
# rows 50-140. This matches a new weekly "Fruta / hortaliza" price report
# entry being recorded ahead of the already-present historical rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 48, shifting rows 48:138 down to 50:140.
$ws.Rows("48:49").Insert()

# --- New row 48: Artic Pride / Primera -------------------------------
$ws.Cells.Item(48, 1).Value  = 8
$ws.Cells.Item(48, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(48, 3).Value  = "Coquimbo"
$ws.Cells.Item(48, 4).Value  = 44536
$ws.Cells.Item(48, 5).Value  = 4
$ws.Cells.Item(48, 6).Value  = "Fruta"
$ws.Cells.Item(48, 7).Value  = 100103
$ws.Cells.Item(48, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(48, 9).Value  = 100103006
$ws.Cells.Item(48, 10).Value = "Nectarín"
$ws.Cells.Item(48, 11).Value = "Artic Pride"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 20
$ws.Cells.Item(48, 14).Value = 460000
$ws.Cells.Item(48, 15).Value = 470000
$ws.Cells.Item(48, 16).Value = 465000
$ws.Cells.Item(48, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(48, 18).Value = "Región Metropolitana"
$ws.Cells.Item(48, 19).Value = 1107
$ws.Cells.Item(48, 20).Value = 420

# --- New row 49: Artic Pride / Segunda --------------------------------
$ws.Cells.Item(49, 1).Value  = 8
$ws.Cells.Item(49, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(49, 3).Value  = "Coquimbo"
$ws.Cells.Item(49, 4).Value  = 44536
$ws.Cells.Item(49, 5).Value  = 4
$ws.Cells.Item(49, 6).Value  = "Fruta"
$ws.Cells.Item(49, 7).Value  = 100103
$ws.Cells.Item(49, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(49, 9).Value  = 100103006
$ws.Cells.Item(49, 10).Value = "Nectarín"
$ws.Cells.Item(49, 11).Value = "Artic Pride"
$ws.Cells.Item(49, 12).Value = "Segunda"
$ws.Cells.Item(49, 13).Value = 16
$ws.Cells.Item(49, 14).Value = 435000
$ws.Cells.Item(49, 15).Value = 440000
$ws.Cells.Item(49, 16).Value = 437500
$ws.Cells.Item(49, 17).Value = "`$/bins (420 kilos)"
$ws.Cells.Item(49, 18).Value = "Región Metropolitana"
$ws.Cells.Item(49, 19).Value = 1042
$ws.Cells.Item(49, 20).Value = 420
